# Auto-generated Excel COM-interop script
# Updates currentAveragePrice / Leve profit columns (H-N) for several
# leve rows across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets,
# reflecting refreshed market-board data from the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1429084.1
$ws.Range("I2").Value = 2500273
$ws.Range("J2").Value = 832.3333
$ws.Range("K2").Value = 2500273
$ws.Range("L2").Value = 832.3333
$ws.Range("M2").Value = -2500160
$ws.Range("N2").Value = -1058.3333

$ws.Range("H64").Value = 7089.222
$ws.Range("I64").Value = 9450
$ws.Range("J64").Value = 5200.6
$ws.Range("K64").Value = 9450
$ws.Range("L64").Value = 5200.6
$ws.Range("M64").Value = -9202
$ws.Range("N64").Value = -5696.6

$ws.Range("H67").Value = 7089.222
$ws.Range("I67").Value = 9450
$ws.Range("J67").Value = 5200.6
$ws.Range("K67").Value = 9450
$ws.Range("L67").Value = 5200.6
$ws.Range("M67").Value = -8592
$ws.Range("N67").Value = -6916.6

$ws.Range("H106").Value = 2580.2
$ws.Range("I106").Value = 2637
$ws.Range("K106").Value = 2637
$ws.Range("M106").Value = -2006

$ws.Range("H128").Value = 43990
$ws.Range("J128").Value = 43990
$ws.Range("L128").Value = 43990
$ws.Range("N128").Value = -53950

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2520.5
$ws.Range("I2").Value = 2304.9285
$ws.Range("J2").Value = 3275
$ws.Range("K2").Value = 2304.9285
$ws.Range("L2").Value = 3275
$ws.Range("M2").Value = -2191.9285
$ws.Range("N2").Value = -3501

$ws.Range("H97").Value = 532.92
$ws.Range("I97").Value = 512.0526
$ws.Range("J97").Value = 599
$ws.Range("K97").Value = 512.0526
$ws.Range("L97").Value = 599
$ws.Range("M97").Value = -16.05259999999998
$ws.Range("N97").Value = -1591

$ws.Range("H116").Value = 2520.5
$ws.Range("I116").Value = 2304.9285
$ws.Range("J116").Value = 3275
$ws.Range("K116").Value = 2304.9285
$ws.Range("L116").Value = 3275
$ws.Range("M116").Value = -10.92849999999999
$ws.Range("N116").Value = -7863

$ws.Range("H132").Value = 1855.7391
$ws.Range("I132").Value = 1159.7858
$ws.Range("J132").Value = 2938.3333
$ws.Range("K132").Value = 3479.3574
$ws.Range("L132").Value = 8814.999899999999
$ws.Range("M132").Value = -949.3574000000003
$ws.Range("N132").Value = -13874.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2520.5
$ws.Range("I3").Value = 2304.9285
$ws.Range("J3").Value = 3275
$ws.Range("K3").Value = 2304.9285
$ws.Range("L3").Value = 3275
$ws.Range("M3").Value = -2190.9285
$ws.Range("N3").Value = -3503

$ws.Range("H94").Value = 744.4286
$ws.Range("I94").Value = 524.6
$ws.Range("J94").Value = 1294
$ws.Range("K94").Value = 524.6
$ws.Range("L94").Value = 1294
$ws.Range("M94").Value = -73.60000000000002
$ws.Range("N94").Value = -2196

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 10021.429
$ws.Range("I62").Value = 10700
$ws.Range("K62").Value = 10700
$ws.Range("M62").Value = -10076

$ws.Range("H65").Value = 10021.429
$ws.Range("I65").Value = 10700
$ws.Range("K65").Value = 53500
$ws.Range("M65").Value = -50380

$ws.Range("H99").Value = 1432216
$ws.Range("I99").Value = 5003506
$ws.Range("J99").Value = 3700
$ws.Range("K99").Value = 5003506
$ws.Range("L99").Value = 3700
$ws.Range("M99").Value = -5002008
$ws.Range("N99").Value = -6696

$ws.Range("H126").Value = 1432216
$ws.Range("I126").Value = 5003506
$ws.Range("J126").Value = 3700
$ws.Range("K126").Value = 15010518
$ws.Range("L126").Value = 11100
$ws.Range("M126").Value = -15008048
$ws.Range("N126").Value = -16040

$ws.Range("H132").Value = 3201.889
$ws.Range("I132").Value = 2165.9
$ws.Range("K132").Value = 6497.700000000001
$ws.Range("M132").Value = -3967.700000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 155.8125
$ws.Range("I33").Value = 165.08333
$ws.Range("J33").Value = 128
$ws.Range("K33").Value = 990.4999799999999
$ws.Range("L33").Value = 768
$ws.Range("M33").Value = -707.4999799999999
$ws.Range("N33").Value = -1334

$ws.Range("H131").Value = 842.05
$ws.Range("I131").Value = 431.1875
$ws.Range("J131").Value = 920.3095
$ws.Range("K131").Value = 1293.5625
$ws.Range("L131").Value = 2760.9285
$ws.Range("M131").Value = 3746.4375
$ws.Range("N131").Value = -12840.9285

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4400.2144
$ws.Range("I70").Value = 4165.615
$ws.Range("K70").Value = 4165.615
$ws.Range("M70").Value = -3895.615

$ws.Range("H73").Value = 4400.2144
$ws.Range("I73").Value = 4165.615
$ws.Range("K73").Value = 4165.615
$ws.Range("M73").Value = -3229.615

$ws.Range("H102").Value = 4504
$ws.Range("I102").Value = 4256
$ws.Range("K102").Value = 4256
$ws.Range("M102").Value = -2634

$ws.Range("H132").Value = 3322.5715
$ws.Range("I132").Value = 2944.889
$ws.Range("K132").Value = 8834.667000000001
$ws.Range("M132").Value = -6304.667000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2158
$ws.Range("I7").Value = 2100.1538
$ws.Range("K7").Value = 2100.1538
$ws.Range("M7").Value = -1988.1538

$ws.Range("H16").Value = 1618.2727
$ws.Range("I16").Value = 1862.625
$ws.Range("J16").Value = 966.6667
$ws.Range("K16").Value = 1862.625
$ws.Range("L16").Value = 966.6667
$ws.Range("M16").Value = -1692.625
$ws.Range("N16").Value = -1306.6667

$ws.Range("H40").Value = 2919.7
$ws.Range("I40").Value = 3042.4285
$ws.Range("K40").Value = 3042.4285
$ws.Range("M40").Value = -2906.4285

$ws.Range("H61").Value = 2599.6924
$ws.Range("I61").Value = 2332.889
$ws.Range("J61").Value = 3200
$ws.Range("K61").Value = 2332.889
$ws.Range("L61").Value = 3200
$ws.Range("M61").Value = -2130.889
$ws.Range("N61").Value = -3604

$ws.Range("H68").Value = 3727.2
$ws.Range("I68").Value = 3810.2856
$ws.Range("K68").Value = 3810.2856
$ws.Range("M68").Value = -3061.2856

$ws.Range("H71").Value = 3727.2
$ws.Range("I71").Value = 3810.2856
$ws.Range("K71").Value = 19051.428
$ws.Range("M71").Value = -15307.428

$ws.Range("H100").Value = 6236290
$ws.Range("I100").Value = 7015563.5
$ws.Range("K100").Value = 7015563.5
$ws.Range("M100").Value = -7015022.5

$ws.Range("H113").Value = 2599.6924
$ws.Range("I113").Value = 2332.889
$ws.Range("J113").Value = 3200
$ws.Range("K113").Value = 2332.889
$ws.Range("L113").Value = 3200
$ws.Range("M113").Value = -162.8890000000001
$ws.Range("N113").Value = -7540

$ws.Range("H122").Value = 4595.1875
$ws.Range("I122").Value = 6186.9473
$ws.Range("K122").Value = 18560.8419
$ws.Range("M122").Value = -16110.8419

$ws.Range("H126").Value = 2158
$ws.Range("I126").Value = 2100.1538
$ws.Range("K126").Value = 6300.4614
$ws.Range("M126").Value = -3830.4614

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7628.5713
$ws.Range("J62").Value = 10000
$ws.Range("L62").Value = 10000
$ws.Range("N62").Value = -11248

$ws.Range("H65").Value = 7628.5713
$ws.Range("J65").Value = 10000
$ws.Range("L65").Value = 50000
$ws.Range("N65").Value = -56240

$ws.Range("H113").Value = 313.5909
$ws.Range("I113").Value = 319.8
$ws.Range("J113").Value = 251.5
$ws.Range("K113").Value = 959.4000000000001
$ws.Range("L113").Value = 754.5
$ws.Range("M113").Value = 1210.6
$ws.Range("N113").Value = -5094.5

$ws.Range("H126").Value = 801.25
$ws.Range("I126").Value = 1000
$ws.Range("J126").Value = 735
$ws.Range("K126").Value = 3000
$ws.Range("L126").Value = 2205
$ws.Range("M126").Value = -530
$ws.Range("N126").Value = -7145

